$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on D/E columns so numeric-looking strings
# (e.g. "211.61", dotted-thousands like "27.557.85") are written back
# as literal text, matching the source inlineStr cells, not converted
# to floating point numbers by Excel's automatic type detection.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '27.557.85'
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -0.58%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.623.89'
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -1.35%  '
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  -0.07%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '211.61'
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  -0.75%  '
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -0.68%  '
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  -0.07%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '23.24'
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  -0.27%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.263'
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  +2.03%  '
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  -0.18%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0889'
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -0.30%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '1.853.24'
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -1.36%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.619.14'
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -1.51%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '4.05'
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  +0.38%  '
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -1.52%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '65.24'
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  +0.80%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '27.521.74'
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '231.69'
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  +0.13%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0719'
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  -0.65%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '7.56'
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -0.88%  '
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -0.04%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '10.45'
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  +3.12%  '
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  +1.18%  '
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  +6.40%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '150.32'
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  -0.41%  '
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -0.47%  '
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  -0.50%  '
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -0.07%  '
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -0.70%  '
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -0.46%  '
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  -0.66%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.467.91'
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  +1.85%  '
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  -1.97%  '
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -2.56%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -0.11%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.952'
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  +7.48%  '
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  +0.80%  '
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  -0.60%  '
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  -2.40%  '
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -0.04%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.02'
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -1.68%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '67.80'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '2.46'
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -1.80%  '
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -1.83%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '5.30'
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -4.78%  '
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  +0.74%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.763.91'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '87.43'
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  +2.21%  '
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  +2.10%  '
